# Slide 4 ("Updates Since Version-02") - Content Placeholder 2
#
# 1) The first bullet ("STAMP Extensions for SR was moved to draft-
#    gandhi-ippm-stamp-srpm", a multi-run/colored paragraph) is removed;
#    the following bullet ("Replaced TWAMP Light draft with STAMP
#    draft") takes its place as a plain, single-run paragraph.
# 2) The bullet "Removed STAMP direct measurement messages" becomes
#    "Removed stand-alone direct measurement messages".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Change 1 ---------------------------------------------------------
# Insert a clean copy of "Replaced TWAMP Light draft with STAMP draft"
# right before the colored multi-run paragraph (so it inherits that
# paragraph's bullet/level formatting plus its first run's plain
# character formatting), then drop both the old multi-run paragraph and
# the paragraph that used to carry this same text lower down.
$stampMoved = $tr.Paragraphs(2, 1)
[void]$stampMoved.InsertBefore("Replaced TWAMP Light draft with STAMP draft`r")
[void]$tr.Paragraphs(3, 1).Delete()   # old "STAMP Extensions for SR..." paragraph
[void]$tr.Paragraphs(3, 1).Delete()   # old duplicate "Replaced TWAMP..." paragraph

# --- Change 2 ---------------------------------------------------------
$removedStamp = $tr.Paragraphs(6, 1)
[void]$removedStamp.InsertBefore("Removed stand-alone direct measurement messages`r")
[void]$tr.Paragraphs(7, 1).Delete()   # old "Removed STAMP direct measurement..." paragraph
